# Fix Training Data Issue (#48)
# The "Date" column (BF) for every team row was off by one day because of
# how the NBA stats site reported the date. Correct each value from the
# old "5-4-2013-14" label to the proper ISO date "2014-05-04".
#
# The new text looks like a date (YYYY-MM-DD), so a plain .Value assignment
# would make Excel auto-convert it into a date serial number. To keep it
# as literal text we temporarily force the range to a text number format
# before writing the values, then clear the (no longer needed) formatting
# so the cells end up back at their original default style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateCol = 58   # column BF
$firstRow = 2
$lastRow = 31
$newDate = "2014-05-04"

$rng = $ws.Range($ws.Cells.Item($firstRow, $dateCol), $ws.Cells.Item($lastRow, $dateCol))

# Force text storage so the ISO-looking string isn't reinterpreted as a date.
$rng.NumberFormat = "@"

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $dateCol).Value = $newDate
}

# Drop the temporary text format again so the cells keep their original
# (default) style, matching the rest of the sheet.
$rng.ClearFormats()
